$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the schedule cells (columns B-F, rows 2-13) with the new assignments.
$ws.Range("B2").Value = "{}"
$ws.Range("C2").Value = "{}"
$ws.Range("D2").Value = "{}"
$ws.Range("E2").Value = "{0: sala nr 10 | Zofia Wiśniewska | Wychowanie fizyczne}"
$ws.Range("F2").Value = "{0: sala nr 3 | Jan Nowak | Język polski}"
$ws.Range("B3").Value = "{}"
$ws.Range("C3").Value = "{}"
$ws.Range("D3").Value = "{}"
$ws.Range("E3").Value = "{0: sala nr 2 | Piotr Wójcik | Biologia}"
$ws.Range("F3").Value = "{0: sala nr 3 | Paweł Lewandowski | Matematyka}"
$ws.Range("B4").Value = "{}"
$ws.Range("C4").Value = "{}"
$ws.Range("D4").Value = "{}"
$ws.Range("E4").Value = "{0: sala nr 5 | Zofia Wiśniewska | Wychowanie fizyczne}"
$ws.Range("F4").Value = "{0: sala nr 11 | Zofia Wiśniewska | Wychowanie fizyczne}"
$ws.Range("B5").Value = "{}"
$ws.Range("C5").Value = "{}"
$ws.Range("D5").Value = "{}"
$ws.Range("E5").Value = "{0: sala nr 3 | Dominik Kaczor | Informatyka}"
$ws.Range("F5").Value = "{0: sala nr 4 | Lena Kowalska | Język angielski}"
$ws.Range("B6").Value = "{}"
$ws.Range("C6").Value = "{}"
$ws.Range("D6").Value = "{0: sala nr 3 | Mateusz Kowalski | Język niemiecki}"
$ws.Range("E6").Value = "{0: sala nr 4 | Piotr Wójcik | Biologia}"
$ws.Range("F6").Value = "{}"
$ws.Range("B7").Value = "{}"
$ws.Range("C7").Value = "{}"
$ws.Range("D7").Value = "{0: sala nr 7 | Karolina Kamińska | Chemia}"
$ws.Range("E7").Value = "{0: sala nr 6 | Mateusz Kowalski | Język niemiecki}"
$ws.Range("F7").Value = "{}"
$ws.Range("B8").Value = "{}"
$ws.Range("C8").Value = "{}"
$ws.Range("D8").Value = "{0: sala nr 6 | Paweł Lewandowski | Matematyka}"
$ws.Range("E8").Value = "{0: sala nr 8 | Dominik Kaczor | Informatyka}"
$ws.Range("F8").Value = "{}"
$ws.Range("B9").Value = "{}"
$ws.Range("C9").Value = "{}"
$ws.Range("D9").Value = "{0: sala nr 1 | Natalia Szymańska | Geografia}"
$ws.Range("E9").Value = "{0: sala nr 6 | Lena Kowalska | Język angielski}"
$ws.Range("F9").Value = "{}"
$ws.Range("B10").Value = "{0: sala nr 5 | Karolina Kamińska | Chemia}"
$ws.Range("C10").Value = "{0: sala nr 3 | Dominik Kaczor | Informatyka}"
$ws.Range("D10").Value = "{0: sala nr 3 | Jan Nowak | Język polski}"
$ws.Range("E10").Value = "{}"
$ws.Range("F10").Value = "{}"
$ws.Range("B11").Value = "{0: sala nr 4 | Paweł Lewandowski | Matematyka}"
$ws.Range("C11").Value = "{0: sala nr 4 | Paweł Lewandowski | Matematyka}"
$ws.Range("D11").Value = "{}"
$ws.Range("E11").Value = "{0: sala nr 3 | Paweł Lewandowski | Matematyka}"
$ws.Range("F11").Value = "{}"
$ws.Range("B12").Value = "{0: sala nr 2 | Katarzyna Mazur | Fizyka}"
$ws.Range("C12").Value = "{0: sala nr 8 | Katarzyna Mazur | Fizyka}"
$ws.Range("D12").Value = "{0: sala nr 3 | Dominik Kaczor | Informatyka}"
$ws.Range("E12").Value = "{0: sala nr 2 | Dominik Kaczor | Informatyka}"
$ws.Range("F12").Value = "{}"
$ws.Range("B13").Value = "{0: sala nr 1 | Jan Nowak | Język polski}"
$ws.Range("C13").Value = "{0: sala nr 9 | Katarzyna Mazur | Fizyka}"
$ws.Range("D13").Value = "{0: sala nr 5 | Natalia Szymańska | Geografia}"
$ws.Range("E13").Value = "{0: sala nr 11 | Lena Kowalska | Język angielski}"
$ws.Range("F13").Value = "{}"

# Update column widths (B, C, D, F) to match the new layout.
# Note: the runtime snaps ColumnWidth to the nearest achievable pixel grid value,
# so we use the closest settable value to the target widths.
$ws.Columns.Item(2).ColumnWidth = 46.833333333333336
$ws.Columns.Item(3).ColumnWidth = 46.833333333333336
$ws.Columns.Item(4).ColumnWidth = 50.833333333333336
$ws.Columns.Item(6).ColumnWidth = 55.833333333333336
